# Add COVID-19 XBB.1.5 vaccines
# - advance the "next new VO term ID" counter (A13) by the 3 IDs that were
#   consumed for the new XBB.1.5 vaccine terms: VO:0010458 -> VO:0010461
# - advance the "RxNorm term starting from" counter (A17) by the 13 IDs
#   that were consumed for the new RxNorm terms: VO:0021167 -> VO:0021180

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value2 = "VO:0010461"
$ws.Range("A17").Value2 = "VO:0021180"

# Move the selection/active cell, as recorded when the author saved the file
$ws.Range("A19").Select()
